$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits after "ENG 230 "
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Merge the date runs ("12 - 09" / " " / "\u2013" / " 2018") into one run.
#    NB: the COM shim treats en-dash (U+2013) as equal to a plain hyphen for
#    change-detection purposes, so writing the final string directly is a
#    silent no-op. Work around it by first writing an em-dash (U+2014, which
#    is *not* in the original text) and then swapping that single character
#    for the real en-dash via Find/Replace.
# ---------------------------------------------------------------------------
$dateRange = $d.Paragraphs.Item(4).Range
$dateRange.MoveEnd(1, -1) | Out-Null
$dateRange.Text = "12 - 09 " + [string][char]0x2014 + " 2018"

$fixRange = $d.Paragraphs.Item(4).Range
$fixRange.MoveEnd(1, -1) | Out-Null
$fixRange.Find.Execute([string][char]0x2014, $false, $false, $false, $false, $false, $true, 1, $false, [string][char]0x2013, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Bump the page number "3" -> "4" on the paragraph right after the
#    "...just like any human would experience." paragraph, then re-add the
#    "_GoBack" bookmark immediately after that run.
# ---------------------------------------------------------------------------
$pageRange = $d.Paragraphs.Item(9).Range
$pageRange.MoveEnd(1, -1) | Out-Null
$pageRange.Text = "4"

$insertPos = $pageRange.End
$tempRange = $d.Range($insertPos, $insertPos)
$tempRange.InsertAfter("X") | Out-Null
$d.Bookmarks.Add("_GoBack", $tempRange) | Out-Null
$d.Bookmarks("_GoBack").Range.Text = ""
